$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the labels in column A:
# A4: "Beitragsbemessungsgrenze PV Ost" -> "Beitragsbemessungsgrenze PV " (trailing space)
# A5: "Beitragsbemessungsgrenze PV West" -> "Jahresarbeitsentgeltgrenze PV " (trailing space)
$ws.Range("A4").Value = "Beitragsbemessungsgrenze PV "
$ws.Range("A5").Value = "Jahresarbeitsentgeltgrenze PV "

# Update the active cell selection to A4
$ws.Range("A4").Select()
